# Tompkins County Health Dept COVID-19 time series - add day 22 (2020-07-09 /
# serial 43928) results. One more positive test result is back; most other
# results are still pending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the bottom of the table (row 23). Since row 22 was the
# last row, Excel carries its formatting down into the freshly inserted row,
# which reproduces the per-column styles (s="5" on B, s="7" on D:F/I:K/M, …).
$ws.Rows.Item(23).Insert()

# Day counter and date: continue the existing "+1 from previous row" pattern.
$ws.Range("A23").Formula = "=A22+1"
$ws.Range("B23").Formula = "=B22+1"

# Raw daily figures.
$ws.Range("C23").Value = 574
$ws.Range("D23").Value = 103
$ws.Range("E23").Value = 1456
$ws.Range("F23").Value = 2133

# These two columns haven't been filled in for this day yet (still "#N/A"
# placeholders), same as every preceding row.
$ws.Range("G23").Value = "#N/A"
$ws.Range("H23").Value = "#N/A"

$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 3

# Recovered/quarantine running total, same running-sum formula as the rows
# above it.
$ws.Range("K23").Formula = "=K22+L23"

$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 69

# Leave the selection on the newly-updated running-total cell, spanning the
# two-row block that now makes up that shared formula.
$excel.Goto($ws.Range("K22:K23"), $true)
